$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.07776
$ws.Range("H2").Value = 3.23328
$ws.Range("I2").Value = 0.1008030812580806
$ws.Range("J2").Value = 0.1008030812580806
$ws.Range("M2").Value = 3.003729333333334
$ws.Range("N2").Value = 9.011188000000001
$ws.Range("O2").Value = 0.2892047983951475
$ws.Range("P2").Value = 0.2892047983951475
$ws.Range("Q2").Value = 3.237299326293334
$ws.Range("R2").Value = 29.13569393664
$ws.Range("S2").Value = 0.02915273479285287
$ws.Range("T2").Value = 0.02915273479285288

# Row 3
$ws.Range("G3").Value = 1.07776
$ws.Range("H3").Value = 3.23328
$ws.Range("I3").Value = 0.1008030812580806
$ws.Range("J3").Value = 0.1008030812580806
$ws.Range("O3").Value = 0.43150426267297
$ws.Range("P3").Value = 0.43150426267297
$ws.Range("Q3").Value = 4.830170407253334
$ws.Range("R3").Value = 43.47153366528001
$ws.Range("S3").Value = 0.04349695925343156
$ws.Range("T3").Value = 0.04349695925343157

# Row 4
$ws.Range("G4").Value = 1.07776
$ws.Range("H4").Value = 3.23328
$ws.Range("I4").Value = 0.1008030812580806
$ws.Range("J4").Value = 0.1008030812580806
$ws.Range("M4").Value = 2.900762333333333
$ws.Range("N4").Value = 8.702287
$ws.Range("O4").Value = 0.2792909389318825
$ws.Range("P4").Value = 0.2792909389318825
$ws.Range("Q4").Value = 3.126325612373333
$ws.Range("R4").Value = 28.13693051136
$ws.Range("S4").Value = 0.02815338721179618
$ws.Range("T4").Value = 0.02815338721179619

# Row 5
$ws.Range("G5").Value = 5.747723
$ws.Range("I5").Value = 0.5375855372420009
$ws.Range("J5").Value = 0.5375855372420009
$ws.Range("M5").Value = 3.003729333333334
$ws.Range("N5").Value = 9.011188000000001
$ws.Range("O5").Value = 0.2892047983951475
$ws.Range("P5").Value = 0.2892047983951475
$ws.Range("Q5").Value = 17.26460417497467
$ws.Range("R5").Value = 155.381437574772
$ws.Range("S5").Value = 0.1554723169182199
$ws.Range("T5").Value = 0.1554723169182199

# Row 6
$ws.Range("G6").Value = 5.747723
$ws.Range("I6").Value = 0.5375855372420009
$ws.Range("J6").Value = 0.5375855372420009
$ws.Range("O6").Value = 0.43150426267297
$ws.Range("P6").Value = 0.43150426267297
$ws.Range("S6").Value = 0.2319704508712621
$ws.Range("T6").Value = 0.2319704508712621

# Row 7
$ws.Range("G7").Value = 5.747723
$ws.Range("I7").Value = 0.5375855372420009
$ws.Range("J7").Value = 0.5375855372420009
$ws.Range("M7").Value = 2.900762333333333
$ws.Range("N7").Value = 8.702287
$ws.Range("O7").Value = 0.2792909389318825
$ws.Range("P7").Value = 0.2792909389318825
$ws.Range("Q7").Value = 16.67277838083367
$ws.Range("R7").Value = 150.055005427503
$ws.Range("S7").Value = 0.1501427694525189
$ws.Range("T7").Value = 0.1501427694525189

# Row 8
$ws.Range("G8").Value = 3.866253666666667
$ws.Range("H8").Value = 11.598761
$ws.Range("I8").Value = 0.3616113814999184
$ws.Range("J8").Value = 0.3616113814999185
$ws.Range("M8").Value = 3.003729333333334
$ws.Range("N8").Value = 9.011188000000001
$ws.Range("O8").Value = 0.2892047983951475
$ws.Range("P8").Value = 0.2892047983951475
$ws.Range("Q8").Value = 11.61317954867422
$ws.Range("R8").Value = 104.518615938068
$ws.Range("S8").Value = 0.1045797466840747
$ws.Range("T8").Value = 0.1045797466840747

# Row 9
$ws.Range("G9").Value = 3.866253666666667
$ws.Range("H9").Value = 11.598761
$ws.Range("I9").Value = 0.3616113814999184
$ws.Range("J9").Value = 0.3616113814999185
$ws.Range("O9").Value = 0.43150426267297
$ws.Range("P9").Value = 0.43150426267297
$ws.Range("Q9").Value = 17.32729369030956
$ws.Range("R9").Value = 155.945643212786
$ws.Range("S9").Value = 0.1560368525482764
$ws.Range("T9").Value = 0.1560368525482764

# Row 10
$ws.Range("G10").Value = 3.866253666666667
$ws.Range("H10").Value = 11.598761
$ws.Range("I10").Value = 0.3616113814999184
$ws.Range("J10").Value = 0.3616113814999185
$ws.Range("M10").Value = 2.900762333333333
$ws.Range("N10").Value = 8.702287
$ws.Range("O10").Value = 0.2792909389318825
$ws.Range("P10").Value = 0.2792909389318825
$ws.Range("Q10").Value = 11.21508300737855
$ws.Range("R10").Value = 100.935747066407
$ws.Range("S10").Value = 0.1009947822675674
$ws.Range("T10").Value = 0.1009947822675674

Write-Host "Updated TPM values for all rows"